$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows after row 11 (before old row 12 "Programa resumido:")
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Row 12: header "Docentes responsáveis:" in column A only
$ws.Cells.Item(12, 1).Value = "Docentes responsáveis:"
$ws.Cells.Item(12, 2).Clear()
$ws.Cells.Item(12, 3).Clear()

# Row 13: teacher 1 name in columns B and C
$ws.Cells.Item(13, 1).Clear()
$ws.Cells.Item(13, 2).Value = "5817344 - Livia Melo Carneiro"
$ws.Cells.Item(13, 3).Value = "5817344 - Livia Melo Carneiro"

# Row 14: teacher 2 name in columns B and C
$ws.Cells.Item(14, 1).Clear()
$ws.Cells.Item(14, 2).Value = "6310296 - Patrícia Caroline Molgero Da Rós"
$ws.Cells.Item(14, 3).Value = "6310296 - Patrícia Caroline Molgero Da Rós"
